# Edit: add "Präparation" body text ("Besonderheiten...") paragraph under
# "Zusatzaufgabe", move the _GoBack bookmark to the end of the new text,
# add the lastRenderedPageBreak hint before the tab run further down, and
# bump the cached PAGE field result in the footer from 2 to 3.

$d = $word.ActiveDocument

# 1) Drop the stale _GoBack bookmark (it will be re-added at the end of the
#    newly inserted paragraph below).
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}

# 2) Fill the first empty "Textkrper" paragraph right after the
#    "Zusatzaufgabe" heading with the new body text (several runs, a couple
#    of proofErr-wrapped technical terms, and the _GoBack bookmark at the
#    very end).
$target = $d.Paragraphs(19)
$target.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Textkrper"/></w:pPr><w:r><w:t>Besonderheiten, die den Regenwurm als Anneliden kennzeichnen sind, dass sein kompletter Körper in Segmente aufgeteilt ist. Ungefähr ab dem 19. Segment wiederholen sich diese in ihrem Inhalt komplett. Er besitzt ein Coelom und bildet damit ein Hydroskelett aus. Weiterhin benutzt er Nephridien um giftige Stoffe auszuscheiden. Seine Außenhaut ist mit Borsten besetzt. Diese helfen bei der Fortbewegung. Er besitzt Sinneszellen zum Unterscheiden</w:t></w:r><w:r><w:t xml:space="preserve">, ob er sich in Licht oder Dunkelheit befindet und Sinneszellen zum Riechen und einen Tastsinn. Diese erfordern größere Knotenpunkte zum Verrechnen der Signale. Deswegen besitzt er ein paariges </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cerebralganglion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Um die Signale über den ganzen Körper zu verteilen hat er ein </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bauchmark</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Er besitzt mehrere Lateralherzen, die das Blut über den Körper verteilen. Die Muskulatur ist schon etwas weiterentwickelt. Er besitzt Längs- und Ringmuskulatur. </w:t></w:r><w:r><w:t xml:space="preserve">Zur Fortpflanzung benutzt er sein </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Clitellum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Alle Regenwürmer sind Zwitter. Wenn der Regenwurm hintere Segmente verliert stirbt er davon nicht. Die Wunde heilt und neue Segmente wachsen heraus. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# 3) Insert the lastRenderedPageBreak marker in front of the tab run of the
#    paragraph with the long signature underline (whole-paragraph rewrite,
#    since InsertXML on a collapsed range is unreliable in this runtime).
$tabsPara = $d.Paragraphs(22)
$tabsPara.Range.InsertXML('<w:p><w:r><w:lastRenderedPageBreak/><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t>____</w:t></w:r><w:r><w:t>_________________________</w:t></w:r><w:r><w:t>______</w:t></w:r></w:p>')

# 4) The document grew by one page: refresh the cached PAGE field result in
#    the footer from "2" to "3".
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2) | Out-Null
